$d = $word.ActiveDocument

$d.Content.Find.Execute("<id>p077r_4</id>", $false, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p077r_4</id>", 2)
$d.Content.Find.Execute("<id>p077v_1</id>", $false, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p077v_1</id>", 2)
$d.Content.Find.Execute("<id>p077v_2</id>", $false, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p077v_2</id>", 2)
$d.Content.Find.Execute("<id>p077v_3</id>", $false, $false, $false, $false, $false,
                         $true, 1, $false, "<id>p077v_3</id>", 2)
